$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the "AUTOMOVEL" field list (column C, rows 16-21) ---
# The old list had an extraneous "ID_MARCA <FK>" entry (marca is reachable
# through modelo, so it doesn't belong directly on automovel) and a
# misplaced "NUM_ANO <FK>" row. Remove the ID_MARCA <FK> row and shift the
# remaining rows up, then drop the now-superfluous trailing row.
$ws.Range("C16").Value = "VAL_PRECO"
$ws.Range("C17").Value = "ID_COR <FK>"
$ws.Range("C18").Value = "DES_CHASSI <PK> "
$ws.Range("C19").Value = "ARQ_FOTO"
$ws.Range("C20:C21").Clear()

# --- Fix over-tagged PK/FK labels on the MARCA / AUTOMOVEL key columns ---
$ws.Range("A27").Value = "NUM_ANO "
$ws.Range("A26").Value = "ID_MARCA <FK> "

# --- Update the saved view position/selection ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E12").Select()
